# Commit: Wed, Jun 10, 2020 11:05:52 AM
#
# 1) Slide 5's table switches from the deck's custom "Table_0" style
#    to PowerPoint's built-in "No Style, No Grid" table style.
# 2) The presentation's theme colour scheme is swapped from the
#    "Integral / Red Violet" palette to the default "Office" palette
#    (font scheme and format scheme are identical between the two
#    themes already, only the 12 theme colours differ).

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 5 -------------------------------------------
$slide5 = $p.Slides.Item(5)
$tableShape = $slide5.Shapes.Item(2)
$tableShape.Table.ApplyStyle("{218F91E5-B8FE-4667-9C2C-155C4A7F86F2}")

# --- 2. Theme colours: Integral (Red Violet) -> Office -------------------
$themeColors = $p.Slides.Item(1).ThemeColorScheme
$themeColors.Item(1).RGB  = 0         # dk1      000000
$themeColors.Item(2).RGB  = 16777215  # lt1      FFFFFF
$themeColors.Item(3).RGB  = 6968388   # dk2      44546A
$themeColors.Item(4).RGB  = 15132391  # lt2      E7E6E6
$themeColors.Item(5).RGB  = 13998939  # accent1  5B9BD5
$themeColors.Item(6).RGB  = 3243501   # accent2  ED7D31
$themeColors.Item(7).RGB  = 10855845  # accent3  A5A5A5
$themeColors.Item(8).RGB  = 49407     # accent4  FFC000
$themeColors.Item(9).RGB  = 12874308  # accent5  4472C4
$themeColors.Item(10).RGB = 4697456   # accent6  70AD47
$themeColors.Item(11).RGB = 12673797  # hlink    0563C1
$themeColors.Item(12).RGB = 7491477   # folHlink 954F72
